# Sicherung vor Änderung SoC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column D (DWPT-Abschnitt?) to 1 for rows 4 through 11
$ws.Range("D4:D11").Value = 1

# Update the active cell selection to D21
$ws.Range("D21").Select()
